$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.028904871456082
$ws.Range("D2").Value = 1.033013221195416
$ws.Range("E2").Value = 1.03818385796169
$ws.Range("F2").Value = 1.047678363721872
$ws.Range("I2").Value = 1.035483074724495
$ws.Range("J2").Value = 1.034054601174966
$ws.Range("K2").Value = 1.035816794069211
$ws.Range("L2").Value = 1.040972615499805
$ws.Range("M2").Value = 1.050440329937776
$ws.Range("N2").Value = 1.015356004987703

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.029707341409009
$ws.Range("D3").Value = 1.033454025046741
$ws.Range("E3").Value = 1.038906953242518
$ws.Range("F3").Value = 1.048506246900083
$ws.Range("I3").Value = 1.035610123831077
$ws.Range("J3").Value = 1.034498630581147
$ws.Range("K3").Value = 1.036067473869659
$ws.Range("L3").Value = 1.041505901525386
$ws.Range("M3").Value = 1.051080064207013
$ws.Range("N3").Value = 1.015504557472633

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.030227332070611
$ws.Range("D4").Value = 1.033739663315654
$ws.Range("E4").Value = 1.039375907060998
$ws.Range("F4").Value = 1.049043166674329
$ws.Range("I4").Value = 1.035691441394329
$ws.Range("J4").Value = 1.034786017460573
$ws.Range("K4").Value = 1.036229357644796
$ws.Range("L4").Value = 1.041851364307769
$ws.Range("M4").Value = 1.051494587953924
$ws.Range("N4").Value = 1.015600664655777

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.030446111262964
$ws.Range("D5").Value = 1.03385984158418
$ws.Range("E5").Value = 1.039573307733649
$ws.Range("F5").Value = 1.049269178678377
$ws.Range("I5").Value = 1.035725413081249
$ws.Range("J5").Value = 1.034906850397865
$ws.Range("K5").Value = 1.036297335086757
$ws.Range("L5").Value = 1.041996688968035
$ws.Range("M5").Value = 1.051668989133637
$ws.Range("N5").Value = 1.015641063755095

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.030482855469182
$ws.Range("D6").Value = 1.033880025614772
$ws.Range("E6").Value = 1.039606466924337
$ws.Range("F6").Value = 1.04930714408772
$ws.Range("I6").Value = 1.035731104493674
$ws.Range("J6").Value = 1.034927139650349
$ws.Range("K6").Value = 1.036308744159924
$ws.Range("L6").Value = 1.042021094975071
$ws.Range("M6").Value = 1.05169827974777
$ws.Range("N6").Value = 1.015647846671965

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.030230254723026
$ws.Range("D7").Value = 1.033741268768988
$ws.Range("E7").Value = 1.039378543748801
$ws.Range("F7").Value = 1.049046185517263
$ws.Range("I7").Value = 1.035691896168323
$ws.Range("J7").Value = 1.034787631976952
$ws.Range("K7").Value = 1.036230266272672
$ws.Range("L7").Value = 1.04185330578298
$ws.Range("M7").Value = 1.051496917781213
$ws.Range("N7").Value = 1.015601204488256

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.029175915856174
$ws.Range("D8").Value = 1.033162107153829
$ws.Range("E8").Value = 1.038428010159133
$ws.Range("F8").Value = 1.047957896020251
$ws.Range("I8").Value = 1.035526195752482
$ws.Range("J8").Value = 1.03420464768024
$ws.Range("K8").Value = 1.035901578594102
$ws.Range("L8").Value = 1.041152760125827
$ws.Range("M8").Value = 1.050656411451289
$ws.Range("N8").Value = 1.015406212062846

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.027323774600915
$ws.Range("D9").Value = 1.032144764573492
$ws.Range("E9").Value = 1.036761274880271
$ws.Range("F9").Value = 1.046049658050443
$ws.Range("I9").Value = 1.035227413873237
$ws.Range("J9").Value = 1.033177949911722
$ws.Range("K9").Value = 1.035319970126967
$ws.Range("L9").Value = 1.039921371279671
$ws.Range("M9").Value = 1.049179797332602
$ws.Range("N9").Value = 1.015062507595404

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.026092978831612
$ws.Range("D10").Value = 1.03146881100693
$ws.Range("E10").Value = 1.035655765850795
$ws.Range("F10").Value = 1.044783990092038
$ws.Range("I10").Value = 1.035023700731247
$ws.Range("J10").Value = 1.032493962770222
$ws.Range("K10").Value = 1.034930677489589
$ws.Range("L10").Value = 1.039102594784243
$ws.Range("M10").Value = 1.048198495304035
$ws.Range("N10").Value = 1.014833330793427

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.025560993700239
$ws.Range("D11").Value = 1.031176679387392
$ws.Range("E11").Value = 1.035178432609213
$ws.Range("F11").Value = 1.044237507999416
$ws.Range("I11").Value = 1.034934426776982
$ws.Range("J11").Value = 1.032197918670539
$ws.Range("K11").Value = 1.034761756380855
$ws.Range("L11").Value = 1.038748584526108
$ws.Range("M11").Value = 1.047774339987372
$ws.Range("N11").Value = 1.014734091177242

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.02536353642167
$ws.Range("D12").Value = 1.031068254919363
$ws.Range("E12").Value = 1.03500133605578
$ws.Range("F12").Value = 1.0440347568316
$ws.Range("I12").Value = 1.034901107229482
$ws.Range("J12").Value = 1.032087975102094
$ws.Range("K12").Value = 1.034698959501569
$ws.Range("L12").Value = 1.038617169842501
$ws.Range("M12").Value = 1.047616905046792
$ws.Range("N12").Value = 1.014697228972916

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.025405885056308
$ws.Range("D13").Value = 1.031091508402371
$ws.Range("E13").Value = 1.035039314501259
$ws.Range("F13").Value = 1.044078236903286
$ws.Range("I13").Value = 1.034908261578771
$ws.Range("J13").Value = 1.032111557426555
$ws.Range("K13").Value = 1.034712431981047
$ws.Range("L13").Value = 1.038645355070946
$ws.Range("M13").Value = 1.047650670132555
$ws.Range("N13").Value = 1.014705136041694

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.025544668837492
$ws.Range("D14").Value = 1.031167715216714
$ws.Range("E14").Value = 1.035163789533416
$ws.Range("F14").Value = 1.044220743685255
$ws.Range("I14").Value = 1.034931675816607
$ws.Range("J14").Value = 1.032188830277177
$ws.Range("K14").Value = 1.034756566625977
$ws.Range("L14").Value = 1.038737720096415
$ws.Range("M14").Value = 1.047761324004965
$ws.Range("N14").Value = 1.014731044137016

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.025630197464907
$ws.Range("D15").Value = 1.031214680234546
$ws.Range("E15").Value = 1.035240510101574
$ws.Range("F15").Value = 1.044308578235371
$ws.Range("I15").Value = 1.034946081024799
$ws.Range("J15").Value = 1.032236443372757
$ws.Range("K15").Value = 1.034783752576464
$ws.Range("L15").Value = 1.038794639949144
$ws.Range("M15").Value = 1.047829516830923
$ws.Range("N15").Value = 1.014747006958829

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.026128305262976
$ws.Range("D16").Value = 1.031488210803519
$ws.Range("E16").Value = 1.035687473700931
$ws.Range("F16").Value = 1.044820291371454
$ws.Range("I16").Value = 1.035029603179151
$ws.Range("J16").Value = 1.032513613028524
$ws.Range("K16").Value = 1.034941880858257
$ws.Range("L16").Value = 1.03912610048107
$ws.Range("M16").Value = 1.048226661155544
$ws.Range("N16").Value = 1.014839916946017

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.026441012776845
$ws.Range("D17").Value = 1.031659940826763
$ws.Range("E17").Value = 1.035968207552951
$ws.Range("F17").Value = 1.045141694819463
$ws.Range("I17").Value = 1.035081709867543
$ws.Range("J17").Value = 1.032687509133164
$ws.Range("K17").Value = 1.035040976438086
$ws.Range("L17").Value = 1.039334158671545
$ws.Range("M17").Value = 1.048475982693059
$ws.Range("N17").Value = 1.014898196055401

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.026623502047176
$ws.Range("D18").Value = 1.031760162014451
$ws.Range("E18").Value = 1.036132085892215
$ws.Range("F18").Value = 1.045329314431808
$ws.Range("I18").Value = 1.035112000027253
$ws.Range("J18").Value = 1.032788951938248
$ws.Range("K18").Value = 1.035098742893802
$ws.Range("L18").Value = 1.039455566079871
$ws.Range("M18").Value = 1.048621480446919
$ws.Range("N18").Value = 1.014932188802114

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.02668574176725
$ws.Range("D19").Value = 1.031794343943546
$ws.Range("E19").Value = 1.036187986334849
$ws.Range("F19").Value = 1.04539331331877
$ws.Range("I19").Value = 1.03512231073113
$ws.Range("J19").Value = 1.03282354333341
$ws.Range("K19").Value = 1.035118433903118
$ws.Range("L19").Value = 1.039496971401498
$ws.Range("M19").Value = 1.048671103700397
$ws.Range("N19").Value = 1.01494377935532

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.026407452663173
$ws.Range("D20").Value = 1.031641510223806
$ws.Range("E20").Value = 1.03593807389953
$ws.Range("F20").Value = 1.045107195697065
$ws.Range("I20").Value = 1.035076129939657
$ws.Range("J20").Value = 1.032668850460587
$ws.Range("K20").Value = 1.035030347961004
$ws.Range("L20").Value = 1.03931183075318
$ws.Range("M20").Value = 1.048449225316909
$ws.Range("N20").Value = 1.014891943303769

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.025503796436726
$ws.Range("D21").Value = 1.031145271813278
$ws.Range("E21").Value = 1.035127129028734
$ws.Range("F21").Value = 1.044178772440753
$ws.Range("I21").Value = 1.034924785296529
$ws.Range("J21").Value = 1.032166074775023
$ws.Range("K21").Value = 1.034743571500109
$ws.Range("L21").Value = 1.038710518664165
$ws.Range("M21").Value = 1.047728736002092
$ws.Range("N21").Value = 1.014723414850051

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.024936475699422
$ws.Range("D22").Value = 1.030833766860028
$ws.Range("E22").Value = 1.034618450191537
$ws.Range("F22").Value = 1.043596405936088
$ws.Range("I22").Value = 1.03482870805254
$ws.Range("J22").Value = 1.031850078385979
$ws.Range("K22").Value = 1.034562963180436
$ws.Range("L22").Value = 1.03833291634058
$ws.Range("M22").Value = 1.047276403351413
$ws.Range("N22").Value = 1.01461745345291

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.025237142564322
$ws.Range("D23").Value = 1.030998853492848
$ws.Range("E23").Value = 1.034887996578766
$ws.Range("F23").Value = 1.043904998850608
$ws.Range("I23").Value = 1.034879727453308
$ws.Range("J23").Value = 1.032017582359473
$ws.Range("K23").Value = 1.034658735142385
$ws.Range("L23").Value = 1.038533045742376
$ws.Range("M23").Value = 1.047516129565426
$ws.Range("N23").Value = 1.014673625512924

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.026422616749104
$ws.Range("D24").Value = 1.031649838053459
$ws.Range("E24").Value = 1.035951689595823
$ws.Range("F24").Value = 1.045122783900998
$ws.Range("I24").Value = 1.035078651586588
$ws.Range("J24").Value = 1.032677281474206
$ws.Range("K24").Value = 1.035035150618511
$ws.Range("L24").Value = 1.039321919622298
$ws.Range("M24").Value = 1.048461315598981
$ws.Range("N24").Value = 1.014894768654607

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.02780190647148
$ws.Range("D25").Value = 1.032407379367934
$ws.Range("E25").Value = 1.03719117898636
$ws.Range("F25").Value = 1.046541848465609
$ws.Range("I25").Value = 1.035305457007136
$ws.Range("J25").Value = 1.033443297345648
$ws.Range("K25").Value = 1.03547060900431
$ws.Range("L25").Value = 1.040239342598692
$ws.Range("M25").Value = 1.04956099747278
$ws.Range("N25").Value = 1.015151372674952

